# Applies the cryptocurrency price/volume refresh described in the commit
# "Updated cryptos list on Thu Nov 30 11:21:00 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell address, new text value, and whether the cell
# needs to be forced to Text format first (otherwise Excel would silently
# reinterpret a numeric-looking string like "228.00" as the number 228).
$changes = @(
    @{ Addr = "D2"; Value = "37.746.62"; ForceText = $false },
    @{ Addr = "E2"; Value = "  -0.89%  "; ForceText = $false },
    @{ Addr = "D3"; Value = "2.033.69"; ForceText = $false },
    @{ Addr = "E3"; Value = "  -0.96%  "; ForceText = $false },
    @{ Addr = "E4"; Value = "  -0.06%  "; ForceText = $false },
    @{ Addr = "D5"; Value = "228.00"; ForceText = $true },
    @{ Addr = "E5"; Value = "  -0.64%  "; ForceText = $false },
    @{ Addr = "D6"; Value = "0.607"; ForceText = $true },
    @{ Addr = "E6"; Value = "  -1.31%  "; ForceText = $false },
    @{ Addr = "D7"; Value = "60.12"; ForceText = $true },
    @{ Addr = "E7"; Value = "  -1.41%  "; ForceText = $false },
    @{ Addr = "E8"; Value = "  -0.03%  "; ForceText = $false },
    @{ Addr = "D9"; Value = "0.376"; ForceText = $true },
    @{ Addr = "E9"; Value = "  -2.35%  "; ForceText = $false },
    @{ Addr = "E10"; Value = "  +2.02%  "; ForceText = $false },
    @{ Addr = "D12"; Value = "14.63"; ForceText = $true },
    @{ Addr = "E12"; Value = "  -1.19%  "; ForceText = $false },
    @{ Addr = "D13"; Value = "2.333.83"; ForceText = $false },
    @{ Addr = "E13"; Value = "  -0.81%  "; ForceText = $false },
    @{ Addr = "D14"; Value = "21.01"; ForceText = $true },
    @{ Addr = "E14"; Value = "  -0.32%  "; ForceText = $false },
    @{ Addr = "D15"; Value = "0.771"; ForceText = $true },
    @{ Addr = "E15"; Value = "  +1.58%  "; ForceText = $false },
    @{ Addr = "D16"; Value = "5.21"; ForceText = $true },
    @{ Addr = "E16"; Value = "  -2.53%  "; ForceText = $false },
    @{ Addr = "D17"; Value = "2.049.54"; ForceText = $false },
    @{ Addr = "E17"; Value = "  -0.47%  "; ForceText = $false },
    @{ Addr = "D18"; Value = "37.730.18"; ForceText = $false },
    @{ Addr = "E18"; Value = "  -0.81%  "; ForceText = $false },
    @{ Addr = "D19"; Value = "69.55"; ForceText = $true },
    @{ Addr = "E19"; Value = "  -0.38%  "; ForceText = $false },
    @{ Addr = "D20"; Value = "5.88"; ForceText = $true },
    @{ Addr = "E20"; Value = "  -6.91%  "; ForceText = $false },
    @{ Addr = "D21"; Value = "0.0₃0823"; ForceText = $false },
    @{ Addr = "E21"; Value = "  -1.05%  "; ForceText = $false },
    @{ Addr = "D22"; Value = "223.59"; ForceText = $true },
    @{ Addr = "E22"; Value = "  -1.12%  "; ForceText = $false },
    @{ Addr = "E23"; Value = "  +0.10%  "; ForceText = $false },
    @{ Addr = "D24"; Value = "2.37"; ForceText = $true },
    @{ Addr = "E24"; Value = "  -2.70%  "; ForceText = $false },
    @{ Addr = "D25"; Value = "2.29"; ForceText = $true },
    @{ Addr = "E25"; Value = "  +2.70%  "; ForceText = $false },
    @{ Addr = "B26"; Value = "Cosmos"; ForceText = $false },
    @{ Addr = "C26"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; ForceText = $false },
    @{ Addr = "D26"; Value = "9.39"; ForceText = $true },
    @{ Addr = "E26"; Value = "  +1.82%  "; ForceText = $false },
    @{ Addr = "B27"; Value = "Monero"; ForceText = $false },
    @{ Addr = "C27"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; ForceText = $false },
    @{ Addr = "D27"; Value = "167.72"; ForceText = $true },
    @{ Addr = "E27"; Value = "  +1.14%  "; ForceText = $false },
    @{ Addr = "E28"; Value = "  -3.35%  "; ForceText = $false },
    @{ Addr = "D29"; Value = "18.76"; ForceText = $true },
    @{ Addr = "E29"; Value = "  -1.45%  "; ForceText = $false },
    @{ Addr = "D30"; Value = "1.27"; ForceText = $true },
    @{ Addr = "E30"; Value = "  -2.34%  "; ForceText = $false },
    @{ Addr = "E31"; Value = "  +0.43%  "; ForceText = $false },
    @{ Addr = "E32"; Value = "  +8.80%  "; ForceText = $false },
    @{ Addr = "D33"; Value = "4.38"; ForceText = $true },
    @{ Addr = "E33"; Value = "  -3.32%  "; ForceText = $false },
    @{ Addr = "D34"; Value = "0.0604"; ForceText = $true },
    @{ Addr = "E34"; Value = "  -0.25%  "; ForceText = $false },
    @{ Addr = "E35"; Value = "  -1.87%  "; ForceText = $false },
    @{ Addr = "D36"; Value = "6.45"; ForceText = $true },
    @{ Addr = "E36"; Value = "  +1.58%  "; ForceText = $false },
    @{ Addr = "D37"; Value = "2.33"; ForceText = $true },
    @{ Addr = "E37"; Value = "  +1.76%  "; ForceText = $false },
    @{ Addr = "D38"; Value = "3.43"; ForceText = $true },
    @{ Addr = "E38"; Value = "  +4.48%  "; ForceText = $false },
    @{ Addr = "E39"; Value = "  +0.11%  "; ForceText = $false },
    @{ Addr = "D40"; Value = "18.01"; ForceText = $true },
    @{ Addr = "E40"; Value = "  +5.53%  "; ForceText = $false },
    @{ Addr = "D41"; Value = "1.537.12"; ForceText = $false },
    @{ Addr = "E41"; Value = "  +1.16%  "; ForceText = $false },
    @{ Addr = "E42"; Value = "  -0.06%  "; ForceText = $false },
    @{ Addr = "D43"; Value = "96.09"; ForceText = $true },
    @{ Addr = "E43"; Value = "  -1.64%  "; ForceText = $false },
    @{ Addr = "E44"; Value = "  -2.10%  "; ForceText = $false },
    @{ Addr = "D45"; Value = "0.0908"; ForceText = $true },
    @{ Addr = "E45"; Value = "  -1.43%  "; ForceText = $false },
    @{ Addr = "D46"; Value = "1.11"; ForceText = $true },
    @{ Addr = "E46"; Value = "  -1.97%  "; ForceText = $false },
    @{ Addr = "D47"; Value = "4.06"; ForceText = $true },
    @{ Addr = "E47"; Value = "  +0.68%  "; ForceText = $false },
    @{ Addr = "E48"; Value = "  -0.73%  "; ForceText = $false },
    @{ Addr = "E49"; Value = "  +0.09%  "; ForceText = $false },
    @{ Addr = "D50"; Value = "7.01"; ForceText = $true },
    @{ Addr = "E50"; Value = "  -0.20%  "; ForceText = $false },
    @{ Addr = "D51"; Value = "2.222.29"; ForceText = $false },
    @{ Addr = "E51"; Value = "  -0.92%  "; ForceText = $false }
)

foreach ($ch in $changes) {
    $cell = $ws.Range($ch.Addr)
    if ($ch.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $ch.Value
}
